$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sponsor rows (name, personType, contact_number, email, previousColab)
$rows = @(
    @{ Row = 8;  Name = "PEPSI";              Type = "J"; Phone = "1234567896"; Email = "pepsi@gmail.com";          Colab = "NO" },
    @{ Row = 9;  Name = "SAMSUNG";            Type = "J"; Phone = "1234567897"; Email = "samsung@hotmail.com";      Colab = "SI" },
    @{ Row = 10; Name = "CARLOS PEREZ";       Type = "N"; Phone = "1234567898"; Email = "carlos@yahoo.com";         Colab = "NO" },
    @{ Row = 11; Name = "ESTEBAN GUTIERREZ";  Type = "N"; Phone = "1234567899"; Email = "esteban@yahoo.com";        Colab = "NO" },
    @{ Row = 12; Name = "MICROSOFT";          Type = "J"; Phone = "1234567800"; Email = "microsoft@outlook.com";    Colab = "SI" }
)

foreach ($r in $rows) {
    $row = $r.Row
    if ($row -eq 11) {
        # This row's data was originally entered phone-number-first, so the
        # shared-string table gained "1234567899" before "ESTEBAN GUTIERREZ".
        $ws.Range("C$row").Value = $r.Phone
        $ws.Range("A$row").Value = $r.Name
    } else {
        $ws.Range("A$row").Value = $r.Name
        $ws.Range("C$row").Value = $r.Phone
    }
    $ws.Range("B$row").Value = $r.Type
    $ws.Range("D$row").Value = $r.Email
    $ws.Range("E$row").Value = $r.Colab
}

# Add mailto hyperlinks for the new email cells, then restore the
# workbook's hyperlink cell style on each (Hyperlinks.Add mints its own
# style the first time it's used, so re-apply the shared one afterwards).
foreach ($r in $rows) {
    $row = $r.Row
    $ws.Hyperlinks.Add($ws.Range("D$row"), "mailto:" + $r.Email)
}
foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Style = "Hipervínculo"
}

[void]$ws.Range("E12").Select()
